$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Progress" column (E)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Progress"
$ws.Range("E2").Value = "1=27"
$ws.Range("E2").NumberFormat = "d-mmm"

# New row 6 - start of a new day entry
$ws.Range("B6").Value = "Lecture 4 Data Mining"
$ws.Range("C6").Value = "10:50 => 12:00"
$ws.Range("D6").Value = "in progress"

# Move selection to B5
$ws.Range("B5").Select()
